$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking price strings are not
# auto-converted to native numbers by Excel (matches original inlineStr cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '29.367.08'
$ws.Range('E2').Value = '  -0.40%  '
$ws.Range('D3').Value = '1.845.59'
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('D4').Value = '0.9983'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '240.33'
$ws.Range('E5').Value = '  -0.55%  '
$ws.Range('D6').Value = '0.6309'
$ws.Range('E6').Value = '  +0.25%  '
$ws.Range('D7').Value = '0.9997'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '0.07534'
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  -0.70%  '
$ws.Range('D10').Value = '24.42'
$ws.Range('E10').Value = '  +0.34%  '
$ws.Range('D11').Value = '0.07700'
$ws.Range('E11').Value = '  +0.03%  '
$ws.Range('D12').Value = '1.861.65'
$ws.Range('E12').Value = '  -1.46%  '
$ws.Range('D13').Value = '4.992'
$ws.Range('E13').Value = '  -0.31%  '
$ws.Range('D14').Value = '0.6829'
$ws.Range('E14').Value = '  -0.71%  '
$ws.Range('D15').Value = '0.000009988'
$ws.Range('E15').Value = '  +1.80%  '
$ws.Range('D16').Value = '82.83'
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('D17').Value = '2.105.18'
$ws.Range('E17').Value = '  -2.75%  '
$ws.Range('D18').Value = '6.116'
$ws.Range('E18').Value = '  -2.17%  '
$ws.Range('D19').Value = '29.390.92'
$ws.Range('D20').Value = '227.75'
$ws.Range('E20').Value = '  -2.62%  '
$ws.Range('D21').Value = '12.43'
$ws.Range('E21').Value = '  -0.51%  '
$ws.Range('D22').Value = '0.9996'
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('D23').Value = '7.538'
$ws.Range('E23').Value = '  -1.25%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').Value = '157.25'
$ws.Range('E25').Value = '  +1.75%  '
$ws.Range('D26').Value = '0.1399'
$ws.Range('E26').Value = '  +0.44%  '
$ws.Range('D27').Value = '8.345'
$ws.Range('E27').Value = '  -1.24%  '
$ws.Range('D28').Value = '17.65'
$ws.Range('E28').Value = '  -0.38%  '
$ws.Range('D29').Value = '1.464'
$ws.Range('E29').Value = '  -0.85%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').Value = '0.05667'
$ws.Range('E30').Value = '  -3.05%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').Value = '1.252'
$ws.Range('E31').Value = '  -0.32%  '
$ws.Range('D32').Value = '4.121'
$ws.Range('E32').Value = '  +0.40%  '
$ws.Range('D33').Value = '4.019'
$ws.Range('E33').Value = '  -0.56%  '
$ws.Range('D34').Value = '1.840'
$ws.Range('E34').Value = '  -2.24%  '
$ws.Range('D35').Value = '1.155'
$ws.Range('E35').Value = '  -1.22%  '
$ws.Range('D36').Value = '0.7117'
$ws.Range('E36').Value = '  -1.03%  '
$ws.Range('D37').Value = '2.593'
$ws.Range('E37').Value = '  +0.24%  '
$ws.Range('D38').Value = '1.260.42'
$ws.Range('E38').Value = '  +1.46%  '
$ws.Range('D39').Value = '0.01812'
$ws.Range('E39').Value = '  +1.55%  '
$ws.Range('D40').Value = '2.776'
$ws.Range('E40').Value = '  -0.68%  '
$ws.Range('D41').Value = '6.227'
$ws.Range('E41').Value = '  +1.13%  '
$ws.Range('D42').Value = '0.9087'
$ws.Range('E42').Value = '  +0.30%  '
$ws.Range('D43').Value = '0.9995'
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('E44').Value = '  -0.88%  '
$ws.Range('D45').Value = '66.18'
$ws.Range('E45').Value = '  -1.57%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '7.060'
$ws.Range('E46').Value = '  -3.49%  '
$ws.Range('B47').Value = 'TheSandbox'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D47').Value = '0.4035'
$ws.Range('E47').Value = '  -0.16%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '9.125'
$ws.Range('E48').Value = '  -0.62%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.00000000117'
$ws.Range('E49').Value = '  -0.95%  '
$ws.Range('E50').Value = '  -1.47%  '
$ws.Range('D51').Value = '0.1125'
$ws.Range('E51').Value = '  +0.63%  '
